{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst last = paras.items[paras.items.length - 1];\nconst r = last.getRange(\"Content\");\nr.load(\"text\");\nawait context.sync();\nlet out = [];\nout.push(\"text=\" + JSON.stringify(r.text));\nr.insertText(r.text, \"Replace\");\nawait context.sync();\nreturn out.join(\"\\n\");\n", "ps1": "$d = $word.ActiveDocument\nfor ($i = 1; $i -le $d.Bookmarks.Count; $i++) {\n  $bm = $d.Bookmarks.Item($i)\n  Write-Output ($i.ToString() + \": \" + $bm.Name)\n}\n"}
